# Update recalculated profit-analysis figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, reflecting
# refreshed market data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 268.9375
$ws.Range("I33").Value = 246.73334
$ws.Range("K33").Value = 246.73334
$ws.Range("M33").Value = -17.73334
$ws.Range("H70").Value = 2951.25
$ws.Range("I70").Value = 1950
$ws.Range("J70").Value = 3285
$ws.Range("K70").Value = 5850
$ws.Range("L70").Value = 9855
$ws.Range("M70").Value = -5580
$ws.Range("N70").Value = -10395
$ws.Range("H73").Value = 2951.25
$ws.Range("I73").Value = 1950
$ws.Range("J73").Value = 3285
$ws.Range("K73").Value = 5850
$ws.Range("L73").Value = 9855
$ws.Range("M73").Value = -4914
$ws.Range("N73").Value = -11727
$ws.Range("H125").Value = 2168
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2168
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 19512
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -24432
$ws.Range("H129").Value = 816.4141
$ws.Range("I129").Value = 297.22223
$ws.Range("J129").Value = 868.3333
$ws.Range("K129").Value = 891.66669
$ws.Range("L129").Value = 2604.9999
$ws.Range("M129").Value = 4108.33331
$ws.Range("N129").Value = -12604.9999
$ws.Range("H135").Value = 866.7083
$ws.Range("I135").Value = 542.5
$ws.Range("K135").Value = 4882.5
$ws.Range("M135").Value = -2347.5
$ws.Range("H140").Value = 83632
$ws.Range("J140").Value = 83632
$ws.Range("L140").Value = 83632
$ws.Range("N140").Value = -93992
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 196.66667
$ws.Range("I4").Value = 195
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 195
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -79
$ws.Range("N4").Value = -432
$ws.Range("H32").Value = 4625.5806
$ws.Range("I32").Value = 2769.7727
$ws.Range("J32").Value = 9162
$ws.Range("K32").Value = 2769.7727
$ws.Range("L32").Value = 9162
$ws.Range("M32").Value = -2482.7727
$ws.Range("N32").Value = -9736
$ws.Range("H61").Value = 1812.1818
$ws.Range("I61").Value = 1536.1538
$ws.Range("J61").Value = 2210.889
$ws.Range("K61").Value = 1536.1538
$ws.Range("L61").Value = 2210.889
$ws.Range("M61").Value = -1324.1538
$ws.Range("N61").Value = -2634.889
$ws.Range("H63").Value = 8660649
$ws.Range("I63").Value = 12594580
$ws.Range("K63").Value = 12594580
$ws.Range("M63").Value = -12593894
$ws.Range("H66").Value = 8660649
$ws.Range("I66").Value = 12594580
$ws.Range("K66").Value = 62972900
$ws.Range("M66").Value = -62969468
$ws.Range("H74").Value = 2183.923
$ws.Range("I74").Value = 784.1429000000001
$ws.Range("J74").Value = 3817
$ws.Range("K74").Value = 784.1429000000001
$ws.Range("L74").Value = 3817
$ws.Range("M74").Value = 89.85709999999995
$ws.Range("N74").Value = -5565
$ws.Range("H76").Value = 27545.5
$ws.Range("J76").Value = 27545.5
$ws.Range("L76").Value = 27545.5
$ws.Range("N76").Value = -28221.5
$ws.Range("H77").Value = 2183.923
$ws.Range("I77").Value = 784.1429000000001
$ws.Range("J77").Value = 3817
$ws.Range("K77").Value = 3920.7145
$ws.Range("L77").Value = 19085
$ws.Range("M77").Value = 447.2855
$ws.Range("N77").Value = -27821
$ws.Range("H79").Value = 27545.5
$ws.Range("J79").Value = 27545.5
$ws.Range("L79").Value = 27545.5
$ws.Range("N79").Value = -29885.5
$ws.Range("H122").Value = 5280.4
$ws.Range("I122").Value = 1170.8
$ws.Range("K122").Value = 3512.4
$ws.Range("M122").Value = -1062.4
$ws.Range("H132").Value = 3324.8572
$ws.Range("I132").Value = 1379
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 4137
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -1607
$ws.Range("N132").Value = -50060
$ws.Range("H136").Value = 1812.1818
$ws.Range("I136").Value = 1536.1538
$ws.Range("J136").Value = 2210.889
$ws.Range("K136").Value = 4608.4614
$ws.Range("L136").Value = 6632.667
$ws.Range("M136").Value = -2058.4614
$ws.Range("N136").Value = -11732.667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 28890
$ws.Range("J118").Value = 28890
$ws.Range("L118").Value = 28890
$ws.Range("N118").Value = -32204
$ws.Range("H134").Value = 1827.7587
$ws.Range("I134").Value = 1357.3214
$ws.Range("K134").Value = 4071.9642
$ws.Range("M134").Value = -1536.9642
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 25626
$ws.Range("J81").Value = 25626
$ws.Range("L81").Value = 25626
$ws.Range("N81").Value = -27622
$ws.Range("H84").Value = 25626
$ws.Range("J84").Value = 25626
$ws.Range("L84").Value = 76878
$ws.Range("N84").Value = -86862
$ws.Range("H109").Value = 34284.5
$ws.Range("J109").Value = 34284.5
$ws.Range("L109").Value = 34284.5
$ws.Range("N109").Value = -36364.5
$ws.Range("H134").Value = 6647.826
$ws.Range("I134").Value = 7329.4116
$ws.Range("K134").Value = 21988.2348
$ws.Range("M134").Value = -19453.2348
$ws.Range("H137").Value = 32429.875
$ws.Range("J137").Value = 32429.875
$ws.Range("L137").Value = 32429.875
$ws.Range("N137").Value = -42629.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 535361.4399999999
$ws.Range("I5").Value = 548.0714
$ws.Range("J5").Value = 1216033
$ws.Range("K5").Value = 1644.2142
$ws.Range("L5").Value = 3648099
$ws.Range("M5").Value = -1532.2142
$ws.Range("N5").Value = -3648323
$ws.Range("H132").Value = 2415.2068
$ws.Range("I132").Value = 830.8461
$ws.Range("J132").Value = 3702.5
$ws.Range("K132").Value = 7477.6149
$ws.Range("L132").Value = 33322.5
$ws.Range("M132").Value = -4947.6149
$ws.Range("N132").Value = -38382.5
$ws.Range("H135").Value = 535361.4399999999
$ws.Range("I135").Value = 548.0714
$ws.Range("J135").Value = 1216033
$ws.Range("K135").Value = 4932.6426
$ws.Range("L135").Value = 10944297
$ws.Range("M135").Value = -2397.6426
$ws.Range("N135").Value = -10949367
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10104568
$ws.Range("I11").Value = 22500000
$ws.Range("J11").Value = 1840947
$ws.Range("K11").Value = 22500000
$ws.Range("L11").Value = 1840947
$ws.Range("M11").Value = -22499861
$ws.Range("N11").Value = -1841225
$ws.Range("H102").Value = 2699.1924
$ws.Range("I102").Value = 1949.8235
$ws.Range("J102").Value = 4114.6665
$ws.Range("K102").Value = 1949.8235
$ws.Range("L102").Value = 4114.6665
$ws.Range("M102").Value = -327.8235
$ws.Range("N102").Value = -7358.6665
$ws.Range("H132").Value = 3415.8076
$ws.Range("I132").Value = 2614.182
$ws.Range("J132").Value = 7824.75
$ws.Range("K132").Value = 7842.545999999999
$ws.Range("L132").Value = 23474.25
$ws.Range("M132").Value = -5312.545999999999
$ws.Range("N132").Value = -28534.25
$ws.Range("H137").Value = 40459.332
$ws.Range("J137").Value = 40459.332
$ws.Range("L137").Value = 40459.332
$ws.Range("N137").Value = -50659.332
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4154.35
$ws.Range("I7").Value = 3819.1333
$ws.Range("J7").Value = 5160
$ws.Range("K7").Value = 3819.1333
$ws.Range("L7").Value = 5160
$ws.Range("M7").Value = -3707.1333
$ws.Range("N7").Value = -5384
$ws.Range("H126").Value = 4154.35
$ws.Range("I126").Value = 3819.1333
$ws.Range("J126").Value = 5160
$ws.Range("K126").Value = 11457.3999
$ws.Range("L126").Value = 15480
$ws.Range("M126").Value = -8987.3999
$ws.Range("N126").Value = -20420
$ws.Range("H136").Value = 4580.5713
$ws.Range("I136").Value = 1324.3334
$ws.Range("K136").Value = 3973.0002
$ws.Range("M136").Value = -1423.0002
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 33910.5
$ws.Range("J93").Value = 33910.5
$ws.Range("L93").Value = 33910.5
$ws.Range("N93").Value = -38902.5
